{"js": "// Adding a note that the instructor should have premade accounts for the students.\n//\n// 1. Insert a new sentence right after \"wiki (Wikipedia.org).\" (and before\n//    \" Each group should:\") in the \"Divide students into groups...\" paragraph.\n// 2. Remove the now-redundant \"Set up an account in the application.\" bullet.\n\n// --- Step 1: insert the new sentence after \"wiki (Wikipedia.org).\" ---\nconst anchorResults = context.document.body.search(\"wiki (Wikipedia.org).\", { matchCase: true });\nanchorResults.load(\"items\");\nawait context.sync();\n\nif (anchorResults.items.length === 0) {\n  throw new Error('Could not find anchor text \"wiki (Wikipedia.org).\"');\n}\n\nanchorResults.items[0].insertText(\n  \" You should provide the students with premade demo accounts so that they do not have to sign up with their own information.\",\n  \"After\"\n);\nawait context.sync();\n\n// --- Step 2: delete the \"Set up an account in the application.\" paragraph ---\nconst targetResults = context.document.body.search(\"Set up an account in the application.\", {\n  matchCase: true\n});\ntargetResults.load(\"items\");\nawait context.sync();\n\nif (targetResults.items.length === 0) {\n  throw new Error('Could not find paragraph \"Set up an account in the application.\"');\n}\n\nconst targetParagraph = targetResults.items[0].paragraphs.getFirst();\ntargetParagraph.delete();\nawait context.sync();\n", "ps1": "# Adding a note that the instructor should have premade accounts for the students.\n#\n# 1. Insert a new sentence right after \"wiki (Wikipedia.org).\" (and before\n#    \" Each group should:\") in the \"Divide students into groups...\" paragraph.\n# 2. Remove the now-redundant \"Set up an account in the application.\" bullet.\n\n$d = $word.ActiveDocument\n\n# --- Step 1: insert the new sentence after \"wiki (Wikipedia.org).\" ---\n$anchor = $d.Content\n$anchorFound = $anchor.Find.Execute(\"wiki (Wikipedia.org).\")\nif (-not $anchorFound) {\n    throw 'Could not find anchor text \"wiki (Wikipedia.org).\"'\n}\n$anchor.InsertAfter(\" You should provide the students with premade demo accounts so that they do not have to sign up with their own information.\")\n\n# --- Step 2: delete the \"Set up an account in the application.\" paragraph ---\n$target = $d.Content\n$targetFound = $target.Find.Execute(\"Set up an account in the application.\")\nif (-not $targetFound) {\n    throw 'Could not find paragraph \"Set up an account in the application.\"'\n}\n$target.Expand(4) | Out-Null  # wdParagraph - grow the range to the whole paragraph, including its mark\n$target.Delete() | Out-Null\n"}
